$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2363.3333
$ws.Range("J17").Value = 2179.4
$ws.Range("L17").Value = 6538.200000000001
$ws.Range("N17").Value = -6874.200000000001
# Row 38
$ws.Range("H38").Value = 3403.25
$ws.Range("J38").Value = 6504
$ws.Range("L38").Value = 19512
$ws.Range("N38").Value = -20256
# Row 41
$ws.Range("H41").Value = 321.5
$ws.Range("I41").Value = 295.33334
$ws.Range("J41").Value = 400
$ws.Range("K41").Value = 295.33334
$ws.Range("L41").Value = 400
$ws.Range("M41").Value = 144.66666
$ws.Range("N41").Value = -1280
# Row 62
$ws.Range("H62").Value = 3932.3333
$ws.Range("J62").Value = 1400
$ws.Range("L62").Value = 1400
$ws.Range("N62").Value = -2648
# Row 65
$ws.Range("H65").Value = 3932.3333
$ws.Range("J65").Value = 1400
$ws.Range("L65").Value = 7000
$ws.Range("N65").Value = -13240
# Row 80
$ws.Range("H80").Value = 636.4
$ws.Range("I80").Value = 710.7143
$ws.Range("K80").Value = 2132.1429
$ws.Range("M80").Value = -1134.1429
# Row 83
$ws.Range("H83").Value = 636.4
$ws.Range("I83").Value = 710.7143
$ws.Range("K83").Value = 6396.428699999999
$ws.Range("M83").Value = -1404.428699999999
# Row 87
$ws.Range("H87").Value = 94569.336
$ws.Range("I87").Value = 93000
$ws.Range("K87").Value = 93000
$ws.Range("M87").Value = -91752
# Row 90
$ws.Range("H90").Value = 94569.336
$ws.Range("I90").Value = 93000
$ws.Range("K90").Value = 279000
$ws.Range("M90").Value = -272760
# Row 131
$ws.Range("H131").Value = 9739.799999999999
$ws.Range("I131").Value = 9739.799999999999
$ws.Range("K131").Value = 29219.4
$ws.Range("M131").Value = -24179.4
# Row 138
$ws.Range("H138").Value = 6422.25
$ws.Range("I138").Value = 2270
$ws.Range("J138").Value = 9651.777
$ws.Range("K138").Value = 6810
$ws.Range("L138").Value = 28955.331
$ws.Range("M138").Value = -1670
$ws.Range("N138").Value = -39235.331

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2995.8333
$ws.Range("I32").Value = 2531.2856
$ws.Range("K32").Value = 2531.2856
$ws.Range("M32").Value = -2244.2856
# Row 56
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
# Row 97
$ws.Range("H97").Value = 2710.8
$ws.Range("I97").Value = 1388.5
$ws.Range("J97").Value = 3592.3333
$ws.Range("K97").Value = 1388.5
$ws.Range("L97").Value = 3592.3333
$ws.Range("M97").Value = -892.5
$ws.Range("N97").Value = -4584.3333
# Row 102
$ws.Range("H102").Value = 1878.25
$ws.Range("I102").Value = 1903.7142
$ws.Range("K102").Value = 1903.7142
$ws.Range("M102").Value = -281.7141999999999
# Row 122
$ws.Range("H122").Value = 5328.2856
$ws.Range("I122").Value = 6041.5
$ws.Range("J122").Value = 1049
$ws.Range("K122").Value = 18124.5
$ws.Range("L122").Value = 3147
$ws.Range("M122").Value = -15674.5
$ws.Range("N122").Value = -8047

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 5
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
# Row 94
$ws.Range("H94").Value = 563.7143
$ws.Range("I94").Value = 407.66666
$ws.Range("J94").Value = 1500
$ws.Range("K94").Value = 407.66666
$ws.Range("L94").Value = 1500
$ws.Range("M94").Value = 43.33334000000002
$ws.Range("N94").Value = -2402
# Row 134
$ws.Range("H134").Value = 3630
$ws.Range("I134").Value = 3630
$ws.Range("K134").Value = 10890
$ws.Range("M134").Value = -8355

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 926.25
$ws.Range("J22").Value = 853.5
$ws.Range("L22").Value = 853.5
$ws.Range("N22").Value = -1553.5
# Row 58
$ws.Range("H58").Value = 3631.5
$ws.Range("I58").Value = 6500
$ws.Range("K58").Value = 6500
$ws.Range("M58").Value = -6297
# Row 132
$ws.Range("H132").Value = 2231.75
$ws.Range("I132").Value = 1348.125
$ws.Range("K132").Value = 4044.375
$ws.Range("M132").Value = -1514.375
# Row 136
$ws.Range("H136").Value = 3631.5
$ws.Range("I136").Value = 6500
$ws.Range("K136").Value = 19500
$ws.Range("M136").Value = -16950

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 104650024
$ws.Range("I4").Value = 64958240
$ws.Range("K4").Value = 194874720
$ws.Range("M4").Value = -194874608

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 16093900
$ws.Range("I11").Value = 16093900
$ws.Range("K11").Value = 16093900
$ws.Range("M11").Value = -16093761
# Row 55
$ws.Range("H55").Value = 9000
$ws.Range("I55").Value = 8000
$ws.Range("J55").Value = 10000
$ws.Range("K55").Value = 8000
$ws.Range("L55").Value = 10000
$ws.Range("M55").Value = -7673
$ws.Range("N55").Value = -10654
# Row 80
$ws.Range("H80").Value = 4947.25
$ws.Range("I80").Value = 2895
$ws.Range("J80").Value = 6999.5
$ws.Range("K80").Value = 2895
$ws.Range("L80").Value = 6999.5
$ws.Range("M80").Value = -1897
$ws.Range("N80").Value = -8995.5
# Row 83
$ws.Range("H83").Value = 4947.25
$ws.Range("I83").Value = 2895
$ws.Range("J83").Value = 6999.5
$ws.Range("K83").Value = 14475
$ws.Range("L83").Value = 34997.5
$ws.Range("M83").Value = -9483
$ws.Range("N83").Value = -44981.5
# Row 102
$ws.Range("H102").Value = 2158.6667
$ws.Range("I102").Value = 1216.7142
$ws.Range("J102").Value = 5455.5
$ws.Range("K102").Value = 1216.7142
$ws.Range("L102").Value = 5455.5
$ws.Range("M102").Value = 405.2858000000001
$ws.Range("N102").Value = -8699.5
# Row 112
$ws.Range("H112").Value = 10293
$ws.Range("J112").Value = 10293
$ws.Range("L112").Value = 10293
$ws.Range("N112").Value = -12509
# Row 122
$ws.Range("H122").Value = 1424
$ws.Range("J122").Value = 1402.3334
$ws.Range("L122").Value = 4207.0002
$ws.Range("N122").Value = -9107.0002
# Row 132
$ws.Range("H132").Value = 3429.7144
$ws.Range("I132").Value = 1505.5
$ws.Range("J132").Value = 4199.4
$ws.Range("K132").Value = 4516.5
$ws.Range("L132").Value = 12598.2
$ws.Range("M132").Value = -1986.5
$ws.Range("N132").Value = -17658.2

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3125.9
$ws.Range("I7").Value = 3084.3333
$ws.Range("K7").Value = 3084.3333
$ws.Range("M7").Value = -2972.3333
# Row 55
$ws.Range("H55").Value = 185.21428
$ws.Range("I55").Value = 176.38461
$ws.Range("J55").Value = 300
$ws.Range("K55").Value = 176.38461
$ws.Range("L55").Value = 300
$ws.Range("M55").Value = -3.384610000000009
$ws.Range("N55").Value = -646
# Row 82
$ws.Range("H82").Value = 1524.5
$ws.Range("I82").Value = 1049.25
$ws.Range("K82").Value = 1049.25
$ws.Range("M82").Value = -688.25
# Row 85
$ws.Range("H85").Value = 1524.5
$ws.Range("I85").Value = 1049.25
$ws.Range("K85").Value = 1049.25
$ws.Range("M85").Value = 198.75
# Row 93
$ws.Range("H93").Value = 2365.182
$ws.Range("I93").Value = 2261.7144
$ws.Range("J93").Value = 2546.25
$ws.Range("K93").Value = 2261.7144
$ws.Range("L93").Value = 2546.25
$ws.Range("M93").Value = -1013.7144
$ws.Range("N93").Value = -5042.25
# Row 126
$ws.Range("H126").Value = 3125.9
$ws.Range("I126").Value = 3084.3333
$ws.Range("K126").Value = 9252.999899999999
$ws.Range("M126").Value = -6782.999899999999
